$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; everything from the old row 20 onward
# shifts down by one (old row 71 becomes new row 72).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record.
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = "2022-05-09"
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 100112035
$ws.Range("G20").Value = "Bruselas (repollito)"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = 30000
$ws.Range("N20").Value = "$/malla 10 kilos"
$ws.Range("O20").Value = "Provincia de Quillota"
$ws.Range("P20").Value = 3000
$ws.Range("Q20").Value = 10
$ws.Range("R20").Value = "Hortaliza"
